# Update test data for the "PC_10" row (row 11) of the product catalog
# test sheet, per "Updated test data of product catalog".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# serial_key: refresh to a newer generated serial
$ws.Range("B11").Value = "307220524Dzl"

# username_pos: swap the placeholder user for the real POS tester account
$ws.Range("E11").Value = "MrunalJagtap"

# new_category_admin: this row is now a "Copy" test case
$ws.Range("S11").Value = "Copy"

# product_name: was the NULL placeholder, now holds a real numeric barcode/id
$ws.Range("U11").Value = 11307101311
